$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B5").Value = "12:25 - 13:30"
$ws.Range("D11").Select()
